$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1733.6
$ws.Range("I19").Value = 1584
$ws.Range("J19").Value = 1833.3334
$ws.Range("K19").Value = 1584
$ws.Range("L19").Value = 1833.3334
$ws.Range("M19").Value = -1409
$ws.Range("N19").Value = -2183.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 77.8
$ws.Range("I33").Value = 79.75
$ws.Range("J33").Value = 70
$ws.Range("K33").Value = 79.75
$ws.Range("L33").Value = 70
$ws.Range("M33").Value = 149.25
$ws.Range("N33").Value = -528

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2119.6
$ws.Range("I80").Value = 1200
$ws.Range("J80").Value = 2513.7144
$ws.Range("K80").Value = 3600
$ws.Range("L80").Value = 7541.1432
$ws.Range("M80").Value = -2602
$ws.Range("N80").Value = -9537.143199999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 2119.6
$ws.Range("I83").Value = 1200
$ws.Range("J83").Value = 2513.7144
$ws.Range("K83").Value = 10800
$ws.Range("L83").Value = 22623.4296
$ws.Range("M83").Value = -5808
$ws.Range("N83").Value = -32607.4296

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 6627.5713
$ws.Range("I86").Value = 6627.5713
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 6627.5713
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -5504.5713

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 6627.5713
$ws.Range("I89").Value = 6627.5713
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 33137.85649999999
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -27521.85649999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 707.41174
$ws.Range("I92").Value = 639.125
$ws.Range("J92").Value = 1800
$ws.Range("K92").Value = 639.125
$ws.Range("L92").Value = 1800
$ws.Range("M92").Value = 608.875
$ws.Range("N92").Value = -4296

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2489.4
$ws.Range("I131").Value = 815.6667
$ws.Range("J131").Value = 5000
$ws.Range("K131").Value = 2447.0001
$ws.Range("L131").Value = 15000
$ws.Range("M131").Value = 2592.9999
$ws.Range("N131").Value = -25080

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4712.8667
$ws.Range("I138").Value = 4196.8
$ws.Range("J138").Value = 4970.9
$ws.Range("K138").Value = 12590.4
$ws.Range("L138").Value = 14912.7
$ws.Range("M138").Value = -7450.400000000001
$ws.Range("N138").Value = -25192.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9900
$ws.Range("I2").Value = 9900
$ws.Range("J2").Value = 9900
$ws.Range("K2").Value = 9900
$ws.Range("L2").Value = 9900
$ws.Range("M2").Value = -9787
$ws.Range("N2").Value = -10126

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 99999
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 99999
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 99999
$ws.Range("N114").Value = -108677

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 9900
$ws.Range("I116").Value = 9900
$ws.Range("J116").Value = 9900
$ws.Range("K116").Value = 9900
$ws.Range("L116").Value = 9900
$ws.Range("M116").Value = -7606
$ws.Range("N116").Value = -14488

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9900
$ws.Range("I3").Value = 9900
$ws.Range("J3").Value = 9900
$ws.Range("K3").Value = 9900
$ws.Range("L3").Value = 9900
$ws.Range("M3").Value = -9786
$ws.Range("N3").Value = -10128

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6180
$ws.Range("I86").Value = 1300
$ws.Range("J86").Value = 9433.333000000001
$ws.Range("K86").Value = 1300
$ws.Range("L86").Value = 9433.333000000001
$ws.Range("M86").Value = -177
$ws.Range("N86").Value = -11679.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 6180
$ws.Range("I89").Value = 1300
$ws.Range("J89").Value = 9433.333000000001
$ws.Range("K89").Value = 6500
$ws.Range("L89").Value = 47166.665
$ws.Range("M89").Value = -884
$ws.Range("N89").Value = -58398.665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3237.4348
$ws.Range("I58").Value = 3165.0476
$ws.Range("J58").Value = 3997.5
$ws.Range("K58").Value = 3165.0476
$ws.Range("L58").Value = 3997.5
$ws.Range("M58").Value = -2962.0476
$ws.Range("N58").Value = -4403.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 1000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 1000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 1000
$ws.Range("N62").Value = -2248

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 1000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 1000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 5000
$ws.Range("N65").Value = -11240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3306.25
$ws.Range("I105").Value = 3081.818
$ws.Range("J105").Value = 3800
$ws.Range("K105").Value = 3081.818
$ws.Range("L105").Value = 3800
$ws.Range("M105").Value = -1334.818
$ws.Range("N105").Value = -7294

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3237.4348
$ws.Range("I136").Value = 3165.0476
$ws.Range("J136").Value = 3997.5
$ws.Range("K136").Value = 9495.1428
$ws.Range("L136").Value = 11992.5
$ws.Range("M136").Value = -6945.1428
$ws.Range("N136").Value = -17092.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 113.333336
$ws.Range("I2").Value = 70
$ws.Range("J2").Value = 135
$ws.Range("K2").Value = 420
$ws.Range("L2").Value = 810
$ws.Range("M2").Value = -307
$ws.Range("N2").Value = -1036

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 969
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 969
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 2907
$ws.Range("N137").Value = -13107
$ws.Range("M137").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 136
$ws.Range("I2").Value = 170.66667
$ws.Range("J2").Value = 66.666664
$ws.Range("K2").Value = 170.66667
$ws.Range("L2").Value = 66.666664
$ws.Range("M2").Value = -57.66667000000001
$ws.Range("N2").Value = -292.666664

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2749
$ws.Range("I80").Value = 2749
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2749
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1751

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2749
$ws.Range("I83").Value = 2749
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 13745
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -8753

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3668.2
$ws.Range("I132").Value = 2397.8
$ws.Range("J132").Value = 4938.6
$ws.Range("K132").Value = 7193.400000000001
$ws.Range("L132").Value = 14815.8
$ws.Range("M132").Value = -4663.400000000001
$ws.Range("N132").Value = -19875.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13969.5
$ws.Range("I22").Value = 23232
$ws.Range("J22").Value = 9999.857
$ws.Range("K22").Value = 23232
$ws.Range("L22").Value = 9999.857
$ws.Range("M22").Value = -22937
$ws.Range("N22").Value = -10589.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 13969.5
$ws.Range("I27").Value = 23232
$ws.Range("J27").Value = 9999.857
$ws.Range("K27").Value = 23232
$ws.Range("L27").Value = 9999.857
$ws.Range("M27").Value = -23125
$ws.Range("N27").Value = -10213.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2624.3
$ws.Range("I40").Value = 2027
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 2027
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -1891
$ws.Range("N40").Value = -8272

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6801420
$ws.Range("I61").Value = 8501167
$ws.Range("J61").Value = 5668255.5
$ws.Range("K61").Value = 8501167
$ws.Range("L61").Value = 5668255.5
$ws.Range("M61").Value = -8500965
$ws.Range("N61").Value = -5668659.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2908.6
$ws.Range("I82").Value = 1444
$ws.Range("J82").Value = 3274.75
$ws.Range("K82").Value = 1444
$ws.Range("L82").Value = 3274.75
$ws.Range("M82").Value = -1083
$ws.Range("N82").Value = -3996.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2908.6
$ws.Range("I85").Value = 1444
$ws.Range("J85").Value = 3274.75
$ws.Range("K85").Value = 1444
$ws.Range("L85").Value = 3274.75
$ws.Range("M85").Value = -196
$ws.Range("N85").Value = -5770.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2399.6
$ws.Range("I100").Value = 2249.5
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2249.5
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1708.5
$ws.Range("N100").Value = -4082

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6801420
$ws.Range("I113").Value = 8501167
$ws.Range("J113").Value = 5668255.5
$ws.Range("K113").Value = 8501167
$ws.Range("L113").Value = 5668255.5
$ws.Range("M113").Value = -8498997
$ws.Range("N113").Value = -5672595.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8808.166999999999
$ws.Range("I62").Value = 4499
$ws.Range("J62").Value = 9670
$ws.Range("K62").Value = 4499
$ws.Range("L62").Value = 9670
$ws.Range("M62").Value = -3875
$ws.Range("N62").Value = -10918

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 8808.166999999999
$ws.Range("I65").Value = 4499
$ws.Range("J65").Value = 9670
$ws.Range("K65").Value = 22495
$ws.Range("L65").Value = 48350
$ws.Range("M65").Value = -19375
$ws.Range("N65").Value = -54590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 998.1
$ws.Range("I107").Value = 998
$ws.Range("J107").Value = 998.3333
$ws.Range("K107").Value = 2994
$ws.Range("L107").Value = 2994.9999
$ws.Range("M107").Value = -1074
$ws.Range("N107").Value = -6834.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = ""
